# [PHOENIX-5879] Modified the BillBasedReceipt according to production dump
#
# The production dump fixed a typo in the "Accountant Category" designation
# (an en-dash was typed instead of a plain hyphen) on the assignmentDetails
# sheet, cell F3. Fixing the text also meant re-keying the cell as Text so
# it won't ever get re-interpreted, and the editor's cursor/selection ended
# up resting on the corrected cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("assignmentDetails")

# Correct the designation text: en-dash "–" -> hyphen "-"
$ws.Range("F3").Value = "Accountant Category - IV"

# Store it as Text (matches production dump's cell format for this column)
$ws.Range("F3").NumberFormat = "@"

# Leave the selection where the edit was made
$ws.Activate() | Out-Null
$ws.Range("F3").Select() | Out-Null

# The workbook window's tab-bar/scrollbar split also moved in the dump
$excel.ActiveWindow.TabRatio = 0.5
